# The "Export" sheet lists account holders sorted by balance (Saldo).
# Three rows (PATRICIA / 005255637, PEDRO / 005081833, DENISE / 004491730)
# need to be removed from the list - they currently sit in sheet rows 7-9
# (row 1 is the header "Conta"/"Nome"/"Saldo").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

$ws.Rows("7:9").Delete()
